$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.639.41"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.972.81"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'244.10"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D7").Value = "'60.04"
$ws.Range("E7").Value = "  +2.42%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.378"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").Value = "'0.0788"
$ws.Range("E10").Value = "  -2.50%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "'14.25"
$ws.Range("E12").Value = "  +4.48%  "
$ws.Range("D13").Value = "'0.844"
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("D14").Value = "2.264.25"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "'21.68"
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("D16").Value = "'5.30"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "1.999.93"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "36.579.24"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "'69.79"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "0.0₃0854"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'229.54"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.08"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +1.62%  "
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("E26").Value = "  +6.38%  "
$ws.Range("D27").Value = "'9.16"
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("D28").Value = "'162.49"
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("D29").Value = "'19.36"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  +18.36%  "
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("D32").Value = "'4.81"
$ws.Range("E32").Value = "  +2.55%  "
$ws.Range("D33").Value = "'0.0615"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").Value = "'4.54"
$ws.Range("E34").Value = "  +7.18%  "
$ws.Range("D35").Value = "'2.28"
$ws.Range("E35").Value = "  +3.19%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "'3.40"
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").Value = "'5.42"
$ws.Range("E39").Value = "  -13.59%  "
$ws.Range("D40").Value = "'0.0967"
$ws.Range("E40").Value = "  -3.65%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").Value = "'15.93"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("D45").Value = "1.368.37"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").Value = "'88.97"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("D48").Value = "'7.23"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").Value = "'46.05"
$ws.Range("E50").Value = "  +5.37%  "
$ws.Range("D51").Value = "2.158.65"
$ws.Range("E51").Value = "  +0.98%  "
